# "student details list problem"
# Highlight the "Student details" bullet item in yellow, matching the
# other highlighted list entries (e.g. Auth, Class, Section, ...).
# This needs to set the highlight on both the run text and the
# paragraph mark, so use the paragraph's own Range (which includes the
# trailing paragraph mark) rather than just the found text range.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Student details") {
        $p.Range.Font.HighlightColorIndex = 7  # wdYellow
    }
}
